# Applies weekly data refresh: rotates the price/volume/date data among
# rows 2, 4, 6, 8, 9 and 10 (rows 3, 5, 7 stay untouched), matching the
# new week's Cebollín price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44225) -> values previously on row 8
$ws.Range("D2").Value = 44210
$ws.Range("J2").Value = 105
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 4000
$ws.Range("M2").Value = 3714
$ws.Range("P2").Value = 1857

# Row 4 (was date 44166) -> values previously on row 10
$ws.Range("D4").Value = 44208
$ws.Range("J4").Value = 85
$ws.Range("K4").Value = 3700
$ws.Range("L4").Value = 4000
$ws.Range("M4").Value = 3824
$ws.Range("N4").Value = "$/paquete 2 kilos"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 1912
$ws.Range("Q4").Value = 2

# Row 6 (was date 44161) -> values previously on row 4
$ws.Range("D6").Value = 44166
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 3500
$ws.Range("L6").Value = 4000
$ws.Range("M6").Value = 3679
$ws.Range("N6").Value = "$/paquete 36 unidades"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 102
$ws.Range("Q6").Value = 36

# Row 8 (was date 44210) -> values previously on row 6
$ws.Range("D8").Value = 44161
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 2800
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 2900
$ws.Range("P8").Value = 1450

# Row 9 (was date 44209) -> values previously on row 2
$ws.Range("D9").Value = 44225
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 3400
$ws.Range("L9").Value = 3700
$ws.Range("M9").Value = 3550
$ws.Range("P9").Value = 1775

# Row 10 (was date 44208) -> values previously on row 9
$ws.Range("D10").Value = 44209
$ws.Range("J10").Value = 150
$ws.Range("K10").Value = 3500
$ws.Range("L10").Value = 4000
$ws.Range("M10").Value = 3767
$ws.Range("P10").Value = 1884
